$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "49.860.67"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.44%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.615.21"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +4.39%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.24%  "

# Row 5
$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "323.69"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.76%  "

# Row 6
$ws.Range("B6").Value = "Solana"
$ws.Range("C6").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "109.97"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.36%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.536"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.75%  "

# Row 8
$ws.Range("E8").Value = "  -0.19%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.564"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.28%  "

# Row 10
$ws.Range("E10").Value = "  +2.69%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.64"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.56%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0827"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.33%  "

# Row 13
$ws.Range("E13").Value = "  +0.67%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.35"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.31%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.017.98"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.04%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.595.69"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.45%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.873"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.41%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "49.787.81"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.56%  "

# Row 19
$ws.Range("E19").Value = "  +13.29%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.39"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.12%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.81"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.65%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0956"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.90%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "282.89"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.95%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "72.85"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.21%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.57"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.45%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.74"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.69%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.01%  "

# Row 28
$ws.Range("E28").Value = "  -1.20%  "

# Row 29
$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.99"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.83%  "

# Row 30
$ws.Range("B30").Value = "Kaspa"
$ws.Range("C30").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.145"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.47%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "36.27"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.36%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.63"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.12%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.78"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.15%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.49"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.59%  "

# Row 35
$ws.Range("E35").Value = "  -0.29%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0796"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.46%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.06"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.13%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.76"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.28%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.07"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.55%  "

# Row 40
$ws.Range("B40").Value = "Stellar"
$ws.Range("C40").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.113"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.05%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "123.49"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.50%  "

# Row 42
$ws.Range("B42").Value = "EnergySwap"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "22.70"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.90%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.23"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.05%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0317"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.97%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.36"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.98%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.049.24"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.35%  "

# Row 47
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.21"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +11.25%  "

# Row 48
$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.04"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +10.79%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.08"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.07%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.39"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.73%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "82.33"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.87%  "
